$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 414, shifting rows 414:539 down to 415:540
$ws.Rows.Item(414).Insert()

# Populate the new row 414 with the data from the edit
$ws.Cells.Item(414, 1).Value = 10
$ws.Cells.Item(414, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(414, 3).Value = "La Araucanía"
$ws.Cells.Item(414, 4).Value = 44876
$ws.Cells.Item(414, 5).Value = 9
$ws.Cells.Item(414, 6).Value = 100112023
$ws.Cells.Item(414, 7).Value = "Brócoli"
$ws.Cells.Item(414, 8).Value = "Sin especificar"
$ws.Cells.Item(414, 9).Value = "Primera"
$ws.Cells.Item(414, 10).Value = 1850
$ws.Cells.Item(414, 11).Value = 800
$ws.Cells.Item(414, 12).Value = 900
$ws.Cells.Item(414, 13).Value = 868
$ws.Cells.Item(414, 14).Value = "$/unidad"
$ws.Cells.Item(414, 15).Value = "Región del Maule"
$ws.Cells.Item(414, 16).Value = 868
$ws.Cells.Item(414, 17).Value = 1
$ws.Cells.Item(414, 18).Value = "Hortaliza"
